# "updated for 30 tarikh"
# Add the 30th (row 33, date 44499) and 31st (row 34, date 44500) meal
# counts, record the two bazar charges for the 29th/30th (F32/F33),
# bump Antor's prior balance (B3), and move the selection to where the
# user was last working (O35). All the SUM()/shared-formula cells
# downstream (B13, U33:U34, F35, K35:U35, K36:T36, K37:T37, B38, ...)
# recompute automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Antor's running balance (B3) went up by 800.
$ws.Range("B3").Value = 5688

# Bazar charges recorded for the 29th and 30th.
$ws.Range("F32").Value = 500
$ws.Range("F33").Value = 300

# Meal counts for the 30th (row 33, 2021-11-30 serial 44499).
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 2
$ws.Range("P33").Value = 2
$ws.Range("Q33").Value = 2
$ws.Range("R33").Value = 2
$ws.Range("S33").Value = 2
$ws.Range("T33").Value = 4

# Meal counts for the 31st (row 34, serial 44500) -- all zero.
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = 0
$ws.Range("S34").Value = 0
$ws.Range("T34").Value = 0

# Leave the cursor where the author last left it.
$ws.Range("O35").Select()
